$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Pages of proceedings" row: convert text values like "S6,941" / "H10,336"
# into plain numeric values with a thousands-separator number format.
$ws.Range("B8").Value = 6941
$ws.Range("C8").Value = 10336

# "Extensions of remarks" row: convert "E1,857" into a plain numeric value.
$ws.Range("C9").Value = 1857

$ws.Range("B8:C8").NumberFormat = "#,##0"
$ws.Range("C9").NumberFormat = "#,##0"

$ws.Range("B8:C8").HorizontalAlignment = -4152
$ws.Range("C9").HorizontalAlignment = -4152
